$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) for rows 2-5: 45184 -> 45185 (i.e. +1 day)
$ws.Range("C2").Value = 45185
$ws.Range("C3").Value = 45185
$ws.Range("C4").Value = 45185
$ws.Range("C5").Value = 45185
